$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "MCT-1A-Circuitos elétricos"

# Row 4
$ws.Range("C4").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("D4").Value = "-"

# Row 6
$ws.Range("C6").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("D6").Value = "MCT-3A-Máquinas Elétricas"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("D7").Value = "MCT-3A-Máquinas Elétricas"

# Row 8
$ws.Range("B8").Value = "-"

# Row 18
$ws.Range("D18").Value = "['ELM-2NA-Automação Industrial', -, -, -]"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("D19").Value = "['ELM-2NA-Automação Industrial', -, -, -]"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("D20").Value = "['ELM-2NA-Automação Industrial', -, -, -]"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("D21").Value = "['ELM-2NA-Automação Industrial', -, -, -]"
$ws.Range("F21").Value = "-"
